# Hotfix: Fri Nov 15 11:48:48 RTZ 2024
# database_tables.xlsx refresh:
#   - Links -> HTML, Releases -> CSS (renamed + reset to placeholder row)
#   - Tasks sheet removed
#   - Python: new Flask file-upload reference row
#   - Bash: clarified pip-install tip + new 'pip install --upgrade' row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename sheets, drop the empty "Tasks" sheet
# ---------------------------------------------------------------
$wb.Worksheets.Item("Links").Name = "HTML"
$wb.Worksheets.Item("Releases").Name = "CSS"
$wb.Worksheets.Item("Tasks").Delete()

# ---------------------------------------------------------------
# 2. Python sheet -> append row 35 (Flask upload snippet)
# ---------------------------------------------------------------
$pyWs = $wb.Worksheets.Item("Python")
$pyB35 = @'
UPLOAD_FOLDER = 'static'
# расширения файлов, которые разрешено загружать
ALLOWED_EXTENSIONS = {'txt', 'pdf', 'png', 'jpg', 'jpeg', 'gif'}
# конфигурируем
app.config['UPLOAD_FOLDER'] = UPLOAD_FOLDER
app.secret_key = "secret key"
@app.route("/upload")
def upload_images():
    return render_template("upload.html")
@app.route('/', methods=['POST'])
def upload_file():
    if request.method == 'POST':
        # check if the post request has the file part
        if 'file' not in request.files:
            flash('No file part')
            return redirect(request.url)
        file = request.files['file']
        if file.filename == '':
            flash('No file selected for uploading')
            return redirect(request.url)
        if file and allowed_file(file.filename):
            filename = secure_filename(file.filename)
            file.save(os.path.join(app.config['UPLOAD_FOLDER'], filename))
            flash('File successfully uploaded')
            return redirect('/upload')
        else:
            flash('Allowed file types are txt, pdf, png, jpg, jpeg, gif, py, docx')
            return redirect(request.url)
<title>Python Flask File Upload Example</title>
<h2>Select a file to upload</h4>
<p>
   {% with messages = get_flashed_messages() %}
     {% if messages %}
      <ul class=flashes>
      {% for message in messages %}
        <li>{{ message }}</li>
      {% endfor %}
      </ul>
     {% endif %}
   {% endwith %}
</p>
<form method="post" action="/" enctype="multipart/form-data">
    <dl>
      <p>
         <input type="file" name="file" autocomplete="off" required>
      </p>
    </dl>
    <p>
      <input type="submit" value="Submit">
   </p>
</form>

'@
$pyC35 = @'
Логика и представление загрузки картинки в приложение
'@
$pyWs.Cells.Item(35, 1).Value = 2111
$pyWs.Cells.Item(35, 2).Value = $pyB35
$pyWs.Cells.Item(35, 3).Value = $pyC35

# ---------------------------------------------------------------
# 3. Bash sheet -> clarify B51, append row 81
# ---------------------------------------------------------------
$bashWs = $wb.Worksheets.Item("Bash")
$bashB51 = @'
Вариант использования с виртуальным окружением:
~/AppData/Local/Programs/Python/Python312/python.exe venv/Scripts/pip.exe install -r requirements.txt
Вариант использования установленным систему PATH:
pip install -r requirements.txt
'@
$bashWs.Cells.Item(51, 2).Value = $bashB51
$bashB81 = @'
pip install --upgrade -r requirements.txt
'@
$bashC81 = @'
Проверка обновлений библиотек
'@
$bashWs.Cells.Item(81, 1).Value = 134
$bashWs.Cells.Item(81, 2).Value = $bashB81
$bashWs.Cells.Item(81, 3).Value = $bashC81

# ---------------------------------------------------------------
# 4. HTML sheet (ex-Links) -> reset to single placeholder "test" row
# ---------------------------------------------------------------
$htmlWs = $wb.Worksheets.Item("HTML")
$htmlWs.Cells.Clear()
$htmlWs.Cells.Item(1, 1).Value = 1
$htmlWs.Cells.Item(1, 2).Value = "test"
$htmlWs.Cells.Item(1, 3).Value = "test"
$htmlWs.Cells.Item(1, 4).Value = "2024-11-13 08:40:02"
$htmlWs.Cells.Item(1, 5).NumberFormat = "General"

# ---------------------------------------------------------------
# 5. CSS sheet (ex-Releases) -> reset to single placeholder "test" row
# ---------------------------------------------------------------
$cssWs = $wb.Worksheets.Item("CSS")
$cssWs.Cells.Clear()
$cssWs.Cells.Item(1, 1).Value = 1
$cssWs.Cells.Item(1, 2).Value = "test"
$cssWs.Cells.Item(1, 3).Value = "test"
$cssWs.Cells.Item(1, 4).Value = "2024-11-13 10:58:04"
$cssWs.Cells.Item(1, 5).NumberFormat = "General"

